$d = $word.ActiveDocument

# Step 1: change " and Java (1 year)." -> " Java (1 year)" (drop leading "and " and trailing ".")
$find = $d.Content.Find
$find.Execute(" and Java (1 year).", $true, $false, $false, $false, $false, $true, 1, $false, " Java (1 year)", 2)

# Step 2: append ",", " and", " SQL (1 year)", "." as separate runs right after "Java (1 year)"
$find2 = $d.Content.Find
$find2.Execute("Java (1 year)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng = $find2.Parent.Duplicate
$rng.Collapse(0)
$rng.InsertAfter(",")

$find3 = $d.Content.Find
$find3.Execute("Java (1 year),", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng3 = $find3.Parent.Duplicate
$rng3.Collapse(0)
$rng3.InsertAfter(" and")

$find4 = $d.Content.Find
$find4.Execute("Java (1 year), and", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng4 = $find4.Parent.Duplicate
$rng4.Collapse(0)
$rng4.InsertAfter(" SQL (1 year)")

$find5 = $d.Content.Find
$find5.Execute("SQL (1 year)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng5 = $find5.Parent.Duplicate
$rng5.Collapse(0)
$rng5.InsertAfter(".")
